$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.014142
$ws.Range("H2").Value = 0.042426
$ws.Range("M2").Value = 3.961421333333333
$ws.Range("N2").Value = 11.884264
$ws.Range("O2").Value = 0.3114993985605504
$ws.Range("P2").Value = 0.3114993985605504
$ws.Range("Q2").Value = 0.056022420496
$ws.Range("R2").Value = 0.504201784464
$ws.Range("S2").Value = 0.3114993985605504
$ws.Range("T2").Value = 0.3114993985605504

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.014142
$ws.Range("H3").Value = 0.042426
$ws.Range("O3").Value = 0.5009735319462221
$ws.Range("P3").Value = 0.500973531946222
$ws.Range("Q3").Value = 0.090098889416
$ws.Range("R3").Value = 0.8108900047439999
$ws.Range("S3").Value = 0.5009735319462221
$ws.Range("T3").Value = 0.500973531946222

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.014142
$ws.Range("H4").Value = 0.042426
$ws.Range("O4").Value = 0.1875270694932276
$ws.Range("P4").Value = 0.1875270694932276
$ws.Range("Q4").Value = 0.033726294144
$ws.Range("R4").Value = 0.303536647296
$ws.Range("S4").Value = 0.1875270694932276
$ws.Range("T4").Value = 0.1875270694932276
